$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-77 shift down to 11-78.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new data record.
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 45063
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = 100112001
$ws.Cells.Item(10, 7).Value = "Berenjena"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 50
$ws.Cells.Item(10, 11).Value = 8500
$ws.Cells.Item(10, 12).Value = 8500
$ws.Cells.Item(10, 13).Value = 8500
$ws.Cells.Item(10, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(10, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 16).Value = 142
$ws.Cells.Item(10, 17).Value = 60
$ws.Cells.Item(10, 18).Value = "Hortaliza"
